# Generate Report for Handback
#
# The d54309a7-... file has now been handed back (in sync with en-US), so its
# row moves to the top of each sheet and gains "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" data. The 2566e582-...
# file is still only "Ready for handoff" and moves down to row 3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Clear every hyperlink on the sheet so we can re-add them in the new order.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ad9315d251f8e6f16756ec152f1367275e4e40e/e2e/d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.md", $null, $null, "d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.md")
$ws.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws.Range("C2").Value2 = "Handed back: in sync with en-US"
$ws.Range("D2").Value2 = "2016-28-17 16:28:33"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ab63c63353cdd86818eede47b11ac320e0a3afb3/e2e/2566e582-7e52-4929-bffc-8bcdd96f74d4.md", $null, $null, "2566e582-7e52-4929-bffc-8bcdd96f74d4.md")
$ws.Range("B3").Value2 = "Ready for handoff"
$ws.Range("C3").Value2 = "Ready for handoff"
$ws.Range("D3").Value2 = "2016-28-17 16:28:15"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A1").Hyperlinks.Delete()

# Row 2: d54309a7 - Handed back, now has target/handback file + datetime.
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ad9315d251f8e6f16756ec152f1367275e4e40e/e2e/d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.md", $null, $null, "d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ad9315d251f8e6f16756ec152f1367275e4e40e/e2e/d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.md", $null, $null, ".md")
$ws.Range("C2").Value2 = "Handed back: in sync with en-US"
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a29e2e76b768f8b6eeb7625a0d637cbf759d1c8d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.a81524c5e0580ecf74bc72631646e917a93ab893.zh-cn.xlf", $null, $null, "d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.a81524c5e0580ecf74bc72631646e917a93ab893.zh-cn.xlf")
$ws.Range("E2").Value2 = "2016-03-17 16:28:30"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ad9315d251f8e6f16756ec152f1367275e4e40e/e2e/d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.md", $null, $null, "d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a29e2e76b768f8b6eeb7625a0d637cbf759d1c8d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.a81524c5e0580ecf74bc72631646e917a93ab893.zh-cn.xlf", $null, $null, "d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.a81524c5e0580ecf74bc72631646e917a93ab893.zh-cn.xlf")
$ws.Range("H2").Value2 = "2016-03-17 16:28:46"
$ws.Range("I2").Value2 = "Include"

# Row 3: 2566e582 - still just Ready for handoff (unchanged data, new row).
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ab63c63353cdd86818eede47b11ac320e0a3afb3/e2e/2566e582-7e52-4929-bffc-8bcdd96f74d4.md", $null, $null, "2566e582-7e52-4929-bffc-8bcdd96f74d4.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/ab63c63353cdd86818eede47b11ac320e0a3afb3/e2e/2566e582-7e52-4929-bffc-8bcdd96f74d4.md", $null, $null, ".md")
$ws.Range("C3").Value2 = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bb4ea0e157a3115d82495b86384de180b388806a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2566e582-7e52-4929-bffc-8bcdd96f74d4.87a998a8cb8eae31c6e9e69739dc92c9921cfdde.zh-cn.xlf", $null, $null, "2566e582-7e52-4929-bffc-8bcdd96f74d4.87a998a8cb8eae31c6e9e69739dc92c9921cfdde.zh-cn.xlf")
$ws.Range("E3").Value2 = "2016-03-17 16:28:12"
$ws.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws.Range("I3").Value2 = "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A1").Hyperlinks.Delete()

# Row 2: d54309a7 - Handed back, now has target/handback file + datetime.
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ad9315d251f8e6f16756ec152f1367275e4e40e/e2e/d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.md", $null, $null, "d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ad9315d251f8e6f16756ec152f1367275e4e40e/e2e/d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.md", $null, $null, ".md")
$ws.Range("C2").Value2 = "Handed back: in sync with en-US"
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a29e2e76b768f8b6eeb7625a0d637cbf759d1c8d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.a81524c5e0580ecf74bc72631646e917a93ab893.de-de.xlf", $null, $null, "d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.a81524c5e0580ecf74bc72631646e917a93ab893.de-de.xlf")
$ws.Range("E2").Value2 = "2016-03-17 16:28:33"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ad9315d251f8e6f16756ec152f1367275e4e40e/e2e/d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.md", $null, $null, "d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7fdbbddb19c32178e751359b1e72d718005e427d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.a81524c5e0580ecf74bc72631646e917a93ab893.de-de.xlf", $null, $null, "d54309a7-b1ef-49b2-a73b-0d7d9956c6ef.a81524c5e0580ecf74bc72631646e917a93ab893.de-de.xlf")
$ws.Range("H2").Value2 = "2016-03-17 16:28:53"
$ws.Range("I2").Value2 = "Include"

# Row 3: 2566e582 - still just Ready for handoff (unchanged data, new row).
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ab63c63353cdd86818eede47b11ac320e0a3afb3/e2e/2566e582-7e52-4929-bffc-8bcdd96f74d4.md", $null, $null, "2566e582-7e52-4929-bffc-8bcdd96f74d4.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/ab63c63353cdd86818eede47b11ac320e0a3afb3/e2e/2566e582-7e52-4929-bffc-8bcdd96f74d4.md", $null, $null, ".md")
$ws.Range("C3").Value2 = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce0b6cf844b815435b69db6a8156669a34c34d5a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2566e582-7e52-4929-bffc-8bcdd96f74d4.87a998a8cb8eae31c6e9e69739dc92c9921cfdde.de-de.xlf", $null, $null, "2566e582-7e52-4929-bffc-8bcdd96f74d4.87a998a8cb8eae31c6e9e69739dc92c9921cfdde.de-de.xlf")
$ws.Range("E3").Value2 = "2016-03-17 16:28:15"
$ws.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws.Range("I3").Value2 = "Include"
